{"js": "// Replace the two-digit multiplication problems in the table with new\n// values, per the commit's diff. Each \"before\" text is unique in the\n// document, so a direct search+replace on each pair is unambiguous.\nconst replacements = [\n  [\"18\u00d761=\", \"96\u00d747=\"],\n  [\"37\u00d753=\", \"13\u00d758=\"],\n  [\"34\u00d718=\", \"81\u00d782=\"],\n  [\"46\u00d712=\", \"95\u00d748=\"],\n  [\"91\u00d758=\", \"80\u00d749=\"],\n  [\"80\u00d787=\", \"30\u00d765=\"],\n  [\"81\u00d750=\", \"86\u00d766=\"],\n  [\"61\u00d744=\", \"33\u00d738=\"],\n  [\"52\u00d741=\", \"14\u00d741=\"],\n  [\"49\u00d758=\", \"16\u00d797=\"],\n  [\"28\u00d772=\", \"35\u00d755=\"],\n  [\"73\u00d749=\", \"95\u00d783=\"],\n  [\"60\u00d740=\", \"71\u00d733=\"],\n  [\"43\u00d765=\", \"25\u00d772=\"],\n  [\"36\u00d728=\", \"74\u00d777=\"],\n  [\"53\u00d719=\", \"70\u00d768=\"],\n  [\"29\u00d780=\", \"32\u00d794=\"],\n  [\"67\u00d754=\", \"54\u00d723=\"],\n  [\"35\u00d771=\", \"43\u00d711=\"],\n  [\"38\u00d778=\", \"80\u00d720=\"],\n  [\"53\u00d780=\", \"25\u00d745=\"],\n  [\"93\u00d727=\", \"12\u00d757=\"],\n  [\"16\u00d791=\", \"72\u00d792=\"],\n  [\"88\u00d725=\", \"72\u00d750=\"],\n  [\"96\u00d795=\", \"96\u00d787=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the table with new\n# values, per the commit's diff. Each \"before\" text is unique in the\n# document, so Find/Replace per pair is unambiguous.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"18\u00d761=\", \"96\u00d747=\"),\n  @(\"37\u00d753=\", \"13\u00d758=\"),\n  @(\"34\u00d718=\", \"81\u00d782=\"),\n  @(\"46\u00d712=\", \"95\u00d748=\"),\n  @(\"91\u00d758=\", \"80\u00d749=\"),\n  @(\"80\u00d787=\", \"30\u00d765=\"),\n  @(\"81\u00d750=\", \"86\u00d766=\"),\n  @(\"61\u00d744=\", \"33\u00d738=\"),\n  @(\"52\u00d741=\", \"14\u00d741=\"),\n  @(\"49\u00d758=\", \"16\u00d797=\"),\n  @(\"28\u00d772=\", \"35\u00d755=\"),\n  @(\"73\u00d749=\", \"95\u00d783=\"),\n  @(\"60\u00d740=\", \"71\u00d733=\"),\n  @(\"43\u00d765=\", \"25\u00d772=\"),\n  @(\"36\u00d728=\", \"74\u00d777=\"),\n  @(\"53\u00d719=\", \"70\u00d768=\"),\n  @(\"29\u00d780=\", \"32\u00d794=\"),\n  @(\"67\u00d754=\", \"54\u00d723=\"),\n  @(\"35\u00d771=\", \"43\u00d711=\"),\n  @(\"38\u00d778=\", \"80\u00d720=\"),\n  @(\"53\u00d780=\", \"25\u00d745=\"),\n  @(\"93\u00d727=\", \"12\u00d757=\"),\n  @(\"16\u00d791=\", \"72\u00d792=\"),\n  @(\"88\u00d725=\", \"72\u00d750=\"),\n  @(\"96\u00d795=\", \"96\u00d787=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair[1], $wdReplaceAll) | Out-Null\n}\n"}
